$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.730.96"
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.469.74"
$ws.Range("E3").Value = "  +2.28%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.81"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.14"
$ws.Range("E6").Value = "  +3.69%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +1.56%  "

$ws.Range("E9").Value = "  -0.99%  "

$ws.Range("E10").Value = "  +1.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("E11").Value = "  +4.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.064.46"
$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.76"
$ws.Range("E13").Value = "  +5.32%  "

$ws.Range("E14").Value = "  +2.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.469.21"
$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.817.10"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.33"
$ws.Range("E18").Value = "  +3.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.37"
$ws.Range("E19").Value = "  +5.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.20"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.74"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.560"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.56"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.611.12"
$ws.Range("E25").Value = "  +2.34%  "

$ws.Range("E26").Value = "  +1.82%  "

$ws.Range("E27").Value = "  -8.42%  "

$ws.Range("E28").Value = "  +2.27%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("E30").Value = "  +2.00%  "

$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("E33").Value = "  -2.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.70"
$ws.Range("E34").Value = "  +1.72%  "

$ws.Range("E35").Value = "  +2.50%  "

$ws.Range("E36").Value = "  +4.58%  "

$ws.Range("E37").Value = "  +7.39%  "

$ws.Range("E38").Value = "  +21.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.08"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("E41").Value = "  -0.34%  "

$ws.Range("E42").Value = "  +2.27%  "

$ws.Range("E43").Value = "  +1.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.23"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +3.72%  "

$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.606.03"
$ws.Range("E47").Value = "  +6.29%  "

$ws.Range("E48").Value = "  +11.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.98"
$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("E50").Value = "  +1.33%  "

$ws.Range("E51").Value = "  -0.10%  "
